$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 435, shifting rows 435:457 down to 436:458
$ws.Rows.Item(435).Insert()

# Populate the newly inserted row 435 with the new data point
$ws.Cells.Item(435, 1).Value = 8
$ws.Cells.Item(435, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(435, 3).Value = "Coquimbo"
$ws.Cells.Item(435, 4).Value = 45147
$ws.Cells.Item(435, 5).Value = 4
$ws.Cells.Item(435, 6).Value = 100112012
$ws.Cells.Item(435, 7).Value = "Espinaca"
$ws.Cells.Item(435, 8).Value = "Sin especificar"
$ws.Cells.Item(435, 9).Value = "Primera"
$ws.Cells.Item(435, 10).Value = 1200
$ws.Cells.Item(435, 11).Value = 450
$ws.Cells.Item(435, 12).Value = 500
$ws.Cells.Item(435, 13).Value = 475
$ws.Cells.Item(435, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(435, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(435, 16).Value = 950
$ws.Cells.Item(435, 17).Value = 0.5
$ws.Cells.Item(435, 18).Value = "Hortaliza"

Write-Output "done"
